# Generate Report for Handoff
# The e5ab4c1b-...md file has been (re-)handed off, which moves its row to the
# bottom of each status table (as the most-recently-touched entry) and
# refreshes its status/timestamps. The other two rows shift up to fill in.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "Overview"
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "ffff24d56b07-c835-4088-9540-3f0a6c1c7f70.md"
$ws.Range("B2").Value = "e2e\ffff24d56b07-c835-4088-9540-3f0a6c1c7f70.md"
$ws.Range("G2").Value = "2016-09-04 15:06:11"

$ws.Range("A3").Value = "ffffff8fe97078-4d50-4430-b6fb-5b6a7aa7a955.md"
$ws.Range("B3").Value = "e2e\ffffff8fe97078-4d50-4430-b6fb-5b6a7aa7a955.md"

$ws.Range("A4").Value = "e5ab4c1b-e5be-460e-bd33-1352fb8b5f31.md"
$ws.Range("B4").Value = "e2e\e5ab4c1b-e5be-460e-bd33-1352fb8b5f31.md"
$ws.Range("E4").Value = "Ready for handoff"
$ws.Range("F4").Value = "Ready for handoff"
$ws.Range("G4").Value = "2016-09-04 15:10:00"

# Rebuild the hyperlinks so their r:id stays attached to the same cell
# position (and therefore the same target URL) while only the display text
# moves with the cell content, matching the reordered rows.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b3d5be1bd97adbba9486cd53e4dba63946ff2693/e2e/e5ab4c1b-e5be-460e-bd33-1352fb8b5f31.md", [Type]::Missing, [Type]::Missing, "e2e\ffff24d56b07-c835-4088-9540-3f0a6c1c7f70.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a4e9c2cfbe37bddaaa9535e2aac374e1a55a59fb/e2e/ffff24d56b07-c835-4088-9540-3f0a6c1c7f70.md", [Type]::Missing, [Type]::Missing, "e2e\ffffff8fe97078-4d50-4430-b6fb-5b6a7aa7a955.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b3d5be1bd97adbba9486cd53e4dba63946ff2693/e2e/ffffff8fe97078-4d50-4430-b6fb-5b6a7aa7a955.md", [Type]::Missing, [Type]::Missing, "e2e\e5ab4c1b-e5be-460e-bd33-1352fb8b5f31.md") | Out-Null

# ------------------------------------------------------------------
# Sheet "zh-cn"
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item(2)

$ws.Range("A2").Value = "ffff24d56b07-c835-4088-9540-3f0a6c1c7f70.md"
$ws.Range("G2").Value = "26ee62c0-6b2f-4dd7-b704-2abfed8096a5.43a82b253913a1b3df6d97b1d24824a26ba7c7e5.zh-cn.xlf"
$ws.Range("H2").Value = "2016-09-04 15:06:05"
$ws.Range("I2").Value = "26ee62c0-6b2f-4dd7-b704-2abfed8096a5.md"
$ws.Range("J2").Value = "26ee62c0-6b2f-4dd7-b704-2abfed8096a5.43a82b253913a1b3df6d97b1d24824a26ba7c7e5.zh-cn.xlf"
$ws.Range("K2").Value = "2016-09-04 15:06:37"

$ws.Range("A3").Value = "ffffff8fe97078-4d50-4430-b6fb-5b6a7aa7a955.md"
$ws.Range("F3").Value = "True"

$ws.Range("A4").Value = "e5ab4c1b-e5be-460e-bd33-1352fb8b5f31.md"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("F4").Value = "False"
$ws.Range("G4").Value = "e5ab4c1b-e5be-460e-bd33-1352fb8b5f31.ffa2d8c6837e2a6232a9e7a0acc7217a08f2a04a.zh-cn.xlf"
$ws.Range("H4").Value = "2016-09-04 15:09:55"
$ws.Range("I4").Value = "e5ab4c1b-e5be-460e-bd33-1352fb8b5f31.md"
$ws.Range("J4").Value = "e5ab4c1b-e5be-460e-bd33-1352fb8b5f31.ffa2d8c6837e2a6232a9e7a0acc7217a08f2a04a.zh-cn.xlf"
$ws.Range("K4").Value = "2016-09-04 15:09:26"
$ws.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b3d5be1bd97adbba9486cd53e4dba63946ff2693/e2e/e5ab4c1b-e5be-460e-bd33-1352fb8b5f31.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/106423395b129f1b11986102662c50a32ec3d70a/e2e/e5ab4c1b-e5be-460e-bd33-1352fb8b5f31.md."

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b3d5be1bd97adbba9486cd53e4dba63946ff2693/e2e/e5ab4c1b-e5be-460e-bd33-1352fb8b5f31.md", [Type]::Missing, [Type]::Missing, "ffff24d56b07-c835-4088-9540-3f0a6c1c7f70.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/4fc36d90acefc1e79d4f7c5a231b5470e9a7fb4f/e2e/e5ab4c1b-e5be-460e-bd33-1352fb8b5f31.md", [Type]::Missing, [Type]::Missing, "26ee62c0-6b2f-4dd7-b704-2abfed8096a5.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a4e9c2cfbe37bddaaa9535e2aac374e1a55a59fb/e2e/ffff24d56b07-c835-4088-9540-3f0a6c1c7f70.md", [Type]::Missing, [Type]::Missing, "ffffff8fe97078-4d50-4430-b6fb-5b6a7aa7a955.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/957355a8d6858a8981897bfa7fec51500f1746fa/e2e/26ee62c0-6b2f-4dd7-b704-2abfed8096a5.md", [Type]::Missing, [Type]::Missing, "26ee62c0-6b2f-4dd7-b704-2abfed8096a5.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b3d5be1bd97adbba9486cd53e4dba63946ff2693/e2e/ffffff8fe97078-4d50-4430-b6fb-5b6a7aa7a955.md", [Type]::Missing, [Type]::Missing, "e5ab4c1b-e5be-460e-bd33-1352fb8b5f31.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/957355a8d6858a8981897bfa7fec51500f1746fa/e2e/26ee62c0-6b2f-4dd7-b704-2abfed8096a5.md", [Type]::Missing, [Type]::Missing, "e5ab4c1b-e5be-460e-bd33-1352fb8b5f31.md") | Out-Null

# Column P (Error Detail) needs to be wide enough for the new long message.
$ws.Columns.Item(16).ColumnWidth = 39.14

# ------------------------------------------------------------------
# Sheet "de-de"
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item(3)

$ws.Range("A2").Value = "ffff24d56b07-c835-4088-9540-3f0a6c1c7f70.md"
$ws.Range("G2").Value = "26ee62c0-6b2f-4dd7-b704-2abfed8096a5.43a82b253913a1b3df6d97b1d24824a26ba7c7e5.de-de.xlf"
$ws.Range("H2").Value = "2016-09-04 15:06:11"
$ws.Range("I2").Value = "26ee62c0-6b2f-4dd7-b704-2abfed8096a5.md"
$ws.Range("J2").Value = "26ee62c0-6b2f-4dd7-b704-2abfed8096a5.43a82b253913a1b3df6d97b1d24824a26ba7c7e5.de-de.xlf"
$ws.Range("K2").Value = "2016-09-04 15:06:44"

$ws.Range("A3").Value = "ffffff8fe97078-4d50-4430-b6fb-5b6a7aa7a955.md"
$ws.Range("F3").Value = "True"

$ws.Range("A4").Value = "e5ab4c1b-e5be-460e-bd33-1352fb8b5f31.md"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("F4").Value = "False"
$ws.Range("G4").Value = "e5ab4c1b-e5be-460e-bd33-1352fb8b5f31.ffa2d8c6837e2a6232a9e7a0acc7217a08f2a04a.de-de.xlf"
$ws.Range("H4").Value = "2016-09-04 15:10:00"
$ws.Range("I4").Value = "e5ab4c1b-e5be-460e-bd33-1352fb8b5f31.md"
$ws.Range("J4").Value = "e5ab4c1b-e5be-460e-bd33-1352fb8b5f31.ffa2d8c6837e2a6232a9e7a0acc7217a08f2a04a.de-de.xlf"
$ws.Range("K4").Value = "2016-09-04 15:09:33"
$ws.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b3d5be1bd97adbba9486cd53e4dba63946ff2693/e2e/e5ab4c1b-e5be-460e-bd33-1352fb8b5f31.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/106423395b129f1b11986102662c50a32ec3d70a/e2e/e5ab4c1b-e5be-460e-bd33-1352fb8b5f31.md."

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b3d5be1bd97adbba9486cd53e4dba63946ff2693/e2e/e5ab4c1b-e5be-460e-bd33-1352fb8b5f31.md", [Type]::Missing, [Type]::Missing, "ffff24d56b07-c835-4088-9540-3f0a6c1c7f70.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/4ea419ec1382221f490891a3b11314d6fa5e67a7/e2e/e5ab4c1b-e5be-460e-bd33-1352fb8b5f31.md", [Type]::Missing, [Type]::Missing, "26ee62c0-6b2f-4dd7-b704-2abfed8096a5.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a4e9c2cfbe37bddaaa9535e2aac374e1a55a59fb/e2e/ffff24d56b07-c835-4088-9540-3f0a6c1c7f70.md", [Type]::Missing, [Type]::Missing, "ffffff8fe97078-4d50-4430-b6fb-5b6a7aa7a955.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/9b4ba8e34870cbe020ee612388a3765851f2179b/e2e/26ee62c0-6b2f-4dd7-b704-2abfed8096a5.md", [Type]::Missing, [Type]::Missing, "26ee62c0-6b2f-4dd7-b704-2abfed8096a5.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b3d5be1bd97adbba9486cd53e4dba63946ff2693/e2e/ffffff8fe97078-4d50-4430-b6fb-5b6a7aa7a955.md", [Type]::Missing, [Type]::Missing, "e5ab4c1b-e5be-460e-bd33-1352fb8b5f31.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/9b4ba8e34870cbe020ee612388a3765851f2179b/e2e/26ee62c0-6b2f-4dd7-b704-2abfed8096a5.md", [Type]::Missing, [Type]::Missing, "e5ab4c1b-e5be-460e-bd33-1352fb8b5f31.md") | Out-Null

# Column P (Error Detail) needs to be wide enough for the new long message.
$ws.Columns.Item(16).ColumnWidth = 39.14

Write-Host "Edit applied"
